# Update cryptos list - prices and 1h volume percentages, and a couple of
# rows whose coin ranking swapped position (37/38 and 44/46).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> @{ col letter = new value }
$updates = @{
    2  = @{ D = "63.579.73"; E = "  +3.06%  " }
    3  = @{ D = "3.489.03";  E = "  +2.27%  " }
    4  = @{ E = "  +0.05%  " }
    5  = @{ D = "581.62";    E = "  +0.71%  " }
    6  = @{ D = "148.09";    E = "  +2.86%  " }
    7  = @{ D = "3.490.06";  E = "  +2.29%  " }
    8  = @{ E = "  -0.04%  " }
    9  = @{ E = "  +0.98%  " }
    10 = @{ E = "  +1.09%  " }
    11 = @{ E = "  +2.30%  " }
    12 = @{ E = "  +5.55%  " }
    13 = @{ D = "4.087.90";  E = "  +2.33%  " }
    14 = @{ D = "29.81";     E = "  +6.47%  " }
    15 = @{ E = "  +2.63%  " }
    16 = @{ D = "3.498.19";  E = "  +2.25%  " }
    17 = @{ D = "0.0000173"; E = "  +1.77%  " }
    18 = @{ D = "63.607.19"; E = "  +2.98%  " }
    19 = @{ D = "6.40";      E = "  +4.29%  " }
    20 = @{ E = "  +5.13%  " }
    21 = @{ D = "9.43";      E = "  +3.09%  " }
    22 = @{ D = "392.07";    E = "  +1.14%  " }
    23 = @{ D = "0.567";     E = "  +3.25%  " }
    24 = @{ D = "75.20";     E = "  +1.16%  " }
    25 = @{ E = "  +0.08%  " }
    26 = @{ D = "3.616.98";  E = "  +1.84%  " }
    27 = @{ D = "0.0000117"; E = "  +1.54%  " }
    28 = @{ E = "  -4.45%  " }
    29 = @{ D = "7.66";      E = "  +3.82%  " }
    30 = @{ D = "1.00";      E = "  +0.10%  " }
    31 = @{ D = "8.29";      E = "  +3.72%  " }
    32 = @{ D = "2.15";      E = "  -0.16%  " }
    33 = @{ E = "  +0.07%  " }
    34 = @{ E = "  -2.54%  " }
    35 = @{ D = "23.70";     E = "  +1.06%  " }
    36 = @{ E = "  +3.49%  " }
    37 = @{ B = "NEARProtocol"; C = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"; D = "5.34";  E = "  +2.98%  " }
    38 = @{ B = "EnergySwap";   C = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens";     D = "32.28"; E = "  +13.93%  " }
    39 = @{ D = "1.60";      E = "  +8.79%  " }
    40 = @{ D = "169.98";    E = "  +0.86%  " }
    41 = @{ D = "3.528.73";  E = "  +2.41%  " }
    42 = @{ D = "0.0772";    E = "  +2.36%  " }
    43 = @{ D = "0.802";     E = "  +2.09%  " }
    44 = @{ B = "Stacks"; C = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"; D = "1.74"; E = "  +4.14%  " }
    45 = @{ E = "  -0.14%  " }
    46 = @{ B = "ONDO"; C = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"; D = "1.22"; E = "  +5.32%  " }
    47 = @{ D = "4.45";      E = "  +0.21%  " }
    48 = @{ D = "2.636.82";  E = "  +5.50%  " }
    49 = @{ E = "  +10.92%  " }
    50 = @{ D = "23.26";     E = "  +2.18%  " }
    51 = @{ D = "6.82";      E = "  +2.92%  " }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
